# Informe Semana 7 - add first thirty scenarios (rows 2-31)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data entry, ordered to reproduce the original shared-string build order ---

# 1) Column G (Realizado Por) for rows 2-31 -> "Javier Estupiñan"
$ws.Range("G2:G31").Value = "Javier Estupiñan"

# 2) B2 - scenario file name for block 1 (rows 2-11)
$ws.Range("B2").Value = "Escenario-prueba_1-10.js"

# 3) C7 - positive-approach description (rows 7-11)
$ws.Range("C7").Value = "Set a-priori, post function, enfoque positivo archivo json"

# 4) C2 - negative-approach description (rows 2-6)
$ws.Range("C2").Value = "Set a-priori, post function, enfoque negativo archivo json"

# 5) Column D (method) for rows 2-31 -> "Post"
$ws.Range("D2:D31").Value = "Post"

# 6) B12 - scenario file name for block 2 (rows 12-21)
$ws.Range("B12").Value = "Escenario-prueba_11-20.js"

# 7) C12 - reuse existing description text (rows 12-21)
$ws.Range("C12").Value = "Set pseudo-aleatorio obtenido de Mockaroo por API, selección aleatoria del dato en el conjunto"

# 8) B22 - scenario file name for block 3 (rows 22-31)
$ws.Range("B22").Value = "Escenario-prueba_21-30.js"

# 9) C22 - reuse existing description text (rows 22-31)
$ws.Range("C22").Value = "Conjunto de datos generada dinámicamente durante la prueba usando faker"

# --- Styling ---

# Column B (rows 2-31): center / top / wrap (new style)
$rngB = $ws.Range("B2:B31")
$rngB.VerticalAlignment = -4160
$rngB.HorizontalAlignment = -4108
$rngB.WrapText = $true

# Column C (rows 2-31): center / center / wrap (existing style already used elsewhere)
$rngC = $ws.Range("C2:C31")
$rngC.VerticalAlignment = -4108
$rngC.HorizontalAlignment = -4108
$rngC.WrapText = $true

# --- Merges for the three ten-row blocks ---
$ws.Range("B2:B11").Merge()
$ws.Range("C2:C6").Merge()
$ws.Range("C7:C11").Merge()
$ws.Range("B12:B21").Merge()
$ws.Range("C12:C21").Merge()
$ws.Range("B22:B31").Merge()
$ws.Range("C22:C31").Merge()

# --- Row height for row 22 ---
$ws.Rows(22).RowHeight = 30

# --- View: freeze header row and scroll/select to reflect the latest edit location ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C22:C31").Select()

# --- Page setup ---
$ws.PageSetup.Orientation = 1
